$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7178
$ws.Range("C3").Value = 159032
$ws.Range("C4").Value = 150088
$ws.Range("C7").Value = 5.62
$ws.Range("C8").Value = 64.14
